$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ivanov, Task 1)
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = '\n     - Catch ValueError for invalid float conversion.\n     - Catch the specific exception you raise on divide-by-zero (ValueError or better use ZeroDivisionError).\n   - Consider letting unexpected exceptions propagate (or at least re-raise after logging) during development.\n\n2. Consistency for divide-by-zero\n   - Python’s built-in division raises ZeroDivisionError. Consider raising ZeroDivisionError in divide(a, b) instead of ValueError to be consistent with standard practice:\n     raise ZeroDivisionError("Cannot divide by zero")\n\n3. Input handling / UX\n   - Trim user input for the operation: operation = input(...).strip()\n   - Consider accepting alternative symbols (e.g., \''x\'' or \''X\'' for multiplication) or validating more explicitly and informing the user.\n   - You may want to loop so the user can perform multiple calculations in one run, and provide a way to quit.\n\n4. Minor robustness notes\n   - Checking b == 0 is fine for direct zeros. If you intend to handle floats that are extremely close to zero, you might consider an epsilon check. For this assignment it’s not necessary, but worth noting.\n   - Consider documenting expected input formats in the module or main docstring.\n\n5. Documentation & style\n   - Add a module-level docstring briefly describing the script.\n   - Consider type hints for function signatures for clarity:\n     def add(a: float, b: float) -> float:\n   - The code is already readable; consider adding a short comment for the main flow or a small README if this were part of a larger project.\n\n6. Tests\n   - Add unit tests (pytest/unittest) for the four operations, including a test for division by zero. This improves correctness assurance and is an excellent habit.\n\nExample small improvements\n- Strip input, catch specific exceptions, and use ZeroDivisionError:\n  operation = input("Enter operation (+, -, *, /): ").strip()\n  try:\n      num1 = float(input("Enter first number: "))\n      num2 = float(input("Enter second number: "))\n  except ValueError:\n      print("Error: Please enter valid numbers.")\n      return\n\n  try:\n      if operation == \''+\'':\n          result = add(num1, num2)\n      ...\n      elif operation == \''/\'':\n          result = divide(num1, num2)\n  except ZeroDivisionError:\n      print("Error: Cannot divide by zero")\n\nFinal summary\n- Functional and well-documented at the function level.\n- A few small improvements around exception specificity, input trimming, minor style touches, and tests would make this excellent.\n- Final overall score: 91 — very good work.'
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 80
$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 45914.6452199074

# Row 3 (Ivanov, Task 2)
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 80
$ws.Range("E3").Value = 'Student and GradeManager classes behave as expected.\n- add_grade validates grade range (0–100).\n- get_average_grade and get_grade_report handle the empty-grades case gracefully.\n- Docstrings for all public classes and methods are provided.\n- Readable and well organized; variable and method names are descriptive.\n\nSuggested improvements and edge cases\n1. Class-average definition\n   - Current implementation computes class average as the mean of each student’s average:\n     total = average(student_i_average) over students.\n     This treats each student equally even if they have different numbers of grades. That is valid in some contexts but often you\''d want the average of all individual grades across the class (i.e., weight each grade equally). Consider documenting which definition you use or provide both options (per-student average vs. per-grade average).\n   - If you intend to exclude students with no grades from the class average, filter them out; currently a student with no grades contributes 0 to the class average which may skew results.\n\n2. Duplicate student IDs\n   - add_student silently overwrites any existing student with the same student_id. It\''s better to:\n     - raise an error (ValueError) if a duplicate ID is added, or\n     - return the existing student and avoid creating a new one, or\n     - provide an update method and document behavior.\n   - At minimum, consider checking and documenting behavior.\n\n3. Type hints and input validation\n   - Add type hints for function signatures (PEP 484) to improve readability and help static analysis tools.\n   - You might validate that subject is a non-empty string and grade is numeric (float/int) before storing.\n\n4. API niceties\n   - Implement __repr__ or __str__ for Student to make debugging and printing nicer (e.g., Student(name=\''Alice\'', id=\''S001\'')).\n   - Consider methods to remove or update grades, list subjects, or export/import student data.\n   - Consider adding unit tests for core behavior and edge cases (no grades, duplicate IDs, invalid grades).\n\n5. Separation of concerns\n   - main() includes interactive/demo code which is fine for a script, but if this module is intended as a reusable library, consider moving demo code to a separate script or protecting main (already done using if __name__ == "__main__":). Also consider making printing optional or providing a CLI/formatter.\n\n6. Small style/format suggestions\n   - When computing averages, you could use statistics.mean for readability.\n   - When there are no students or grades, decide whether to return None (to indicate absence) rather than 0, depending on the API contract; returning 0 is okay but document it clearly.\n   - Add unit tests or a README describing intended behavior (how class average is computed, how duplicates are handled, etc.).\n\nExample changes (brief)\n- Exclude students with no grades from class average:\n  students_with_grades = [s for s in self.students.values() if s.grades]\n  if not students_with_grades:\n      return 0\n  total_average = sum(s.get_average_grade() for s in students_with_grades)\n  return total_average / len(students_with_grades)\n\n- Prevent overwriting student IDs:\n  if student_id in self.students:\n      raise ValueError(f"Student with ID {student_id} already exists.")\n  student = Student(name, student_id)\n  self.students[student_id] = student\n\nSummary\n- The submission meets the assignment and is implemented cleanly.\n- Minor improvements would make the code more robust and production-ready: clarify class-average definition, handle duplicate IDs, add type hints and more validation, and possibly add __repr__ and unit tests.\n\nWell done — solid implementation with clear scope.'
$ws.Range("F3").Value = 80
$ws.Range("G3").Value = 80
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 45914.6452199074

# Row 4 (Petrov, Task 1) - only the Task column type changes from text to number
$ws.Range("C4").Value = 1

# Row 5 (Petrov, Task 2) - only the Task column type changes from text to number
$ws.Range("C5").Value = 2
